$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1198
$ws.Range("I2").Value = 3116
$ws.Range("J2").Value = 13144
$ws.Range("K2").Value = 72
$ws.Range("L2").Value = 3582
$ws.Range("M2").Value = 237
$ws.Range("N2").Value = 2270
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 53
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = 192
$ws.Range("S2").Value = 1409
$ws.Range("T2").Value = 2281
$ws.Range("U2").Value = 186
$ws.Range("V2").Value = 20378
$ws.Range("W2").Value = 9
$ws.Range("X2").Value = 20329
$ws.Range("Y2").Value = 34
$ws.Range("Z2").Value = 314
$ws.Range("AA2").Value = 149
